# Fix the "Treatment Agent" SQL expression in the TreatmentTab query (cell B5):
# remove the redundant CONCAT(...) wrapper around REPLACE(...).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldFragment = "CONCAT(REPLACE(trt.treatment_agent, ';', ', '))"
$newFragment = "REPLACE(trt.treatment_agent, ';', ', ')"

$cell = $ws.Range("B5")
$cell.Value = $cell.Value2.Replace($oldFragment, $newFragment)

# Reflect the user's subsequent click on C5 (selection moved from C4 to C5,
# and the view scrolled back to the top of the sheet).
$ws.Range("C5").Select()
